$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.186905860900879
$ws.Range("B1").Value = 2.179807901382446
$ws.Range("C1").Value = 3.686401128768921
$ws.Range("D1").Value = 3.305572032928467
$ws.Range("E1").Value = 1.143163442611694
